# This script updates the "Requirements" worksheet of the AQUATOX
# Models and Requirements workbook to reflect completion of the
# Organic Matter model's data-requirements checking work.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Requirements")

# Swap B34/B35 so the "Six OM State Variables" requirement lines up with
# the AQTVolumeModel requirement row the way the corrected sheet orders them.
$ws.Range("B34").Value = "AQTVolumeModel"
$ws.Range("B35").Value = "Six OM State Variables"

# The Organic Matter model's CheckDataRequirements routine has now been
# written, so replace the "TBA" placeholders in column C (rows 34-37) with
# the real test-script reference.
$ws.Range("C34").Value = "AQTOrganicMatter.CheckDataRequirements"
$ws.Range("C35").Value = "AQTOrganicMatter.CheckDataRequirements"
$ws.Range("C36").Value = "AQTOrganicMatter.CheckDataRequirements"
$ws.Range("C37").Value = "AQTOrganicMatter.CheckDataRequirements"

# Fix the typo "depositin" -> "deposition" in the erosion/deposition
# soft-requirement row.
$ws.Range("B39").Value = "Erosion / deposition rates, soft requirement"

# Move the sheet's selection to reflect where the author was last working.
$ws.Range("A41").Select()
